$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "Asignacion" worksheet between "catalogo" and "Hoja2"
# ------------------------------------------------------------------
$catalogo = $wb.Worksheets.Item("catalogo")
$hoja2    = $wb.Worksheets.Item("Hoja2")

$asig = $wb.Worksheets.Add($null, $catalogo)
$asig.Name = "Asignacion"

# ------------------------------------------------------------------
# 2. Header row (row 3): "Funcionalidad" / "Responsable" (merged C3:D3)
#    Reuse the existing bold header look (font+fill+border) from the
#    "catalogo" sheet, then bump the font size to 14 and center it.
# ------------------------------------------------------------------
$catalogo.Range("D7").Copy()
$asig.Range("B3").PasteSpecial(-4122)
$asig.Range("B3").Font.Size = 14
$asig.Range("B3").HorizontalAlignment = -4108

$catalogo.Range("D7").Copy()
$asig.Range("C3:D3").PasteSpecial(-4122)
$asig.Range("C3:D3").Font.Size = 14
$asig.Range("C3:D3").HorizontalAlignment = -4108

$asig.Range("B3").Value = "Funcionalidad"
$asig.Range("C3").Value = "Responsable"
$asig.Range("C3:D3").Merge()
$asig.Rows.Item(3).RowHeight = 19

# ------------------------------------------------------------------
# 3. Numbered list (B4:B12) + Funcionalidad/Responsable pairs (C/D4:12)
#    Reuse the plain bordered look already used on "catalogo".
# ------------------------------------------------------------------
$catalogo.Range("D8").Copy()
$asig.Range("B4:B12").PasteSpecial(-4122)
$asig.Range("B4:B12").HorizontalAlignment = -4108

$catalogo.Range("D8").Copy()
$asig.Range("C4:D12").PasteSpecial(-4122)

$asig.Range("B4").Value = 1
$asig.Range("B5").Value = 2
$asig.Range("B6").Value = 3
$asig.Range("B7").Value = 4
$asig.Range("B8").Value = 5
$asig.Range("B9").Value = 6
$asig.Range("B10").Value = 7
$asig.Range("B11").Value = 8
$asig.Range("B12").Value = 9

$asig.Range("C4").Value = "JHONATAN STEVEN"
$asig.Range("C5").Value = "MIGUEL ANGEL"
$asig.Range("D5").Value = "JHOHANNS"
$asig.Range("C6").Value = "JAIRO ANDRES"
$asig.Range("C7").Value = "BRIAN DAVID"
$asig.Range("C8").Value = "BRAYAN ESTIVEN"
$asig.Range("C9").Value = "SERGIO ARMANDO"
$asig.Range("C10").Value = "LAURA MILENA"
$asig.Range("C11").Value = "JUAN SEBASTIAN"
$asig.Range("C12").Value = "ALEJANDRO"

# ------------------------------------------------------------------
# 4. Hidden helper lookup table (G4:H13) used by the sheet
# ------------------------------------------------------------------
$asig.Range("G4").Value = 1
$asig.Range("G5").Value = 2
$asig.Range("G6").Value = 3
$asig.Range("G7").Value = 4
$asig.Range("G8").Value = 5
$asig.Range("G9").Value = 6
$asig.Range("G10").Value = 7
$asig.Range("G11").Value = 8
$asig.Range("G12").Value = 9
$asig.Range("G13").Value = 10

$catalogo.Range("D8").Copy()
$asig.Range("H4:H13").PasteSpecial(-4122)

$asig.Range("H4").Value = "JUAN SEBASTIAN"
$asig.Range("H5").Value = "SERGIO ARMANDO"
$asig.Range("H6").Value = "BRAYAN ESTIVEN"
$asig.Range("H7").Value = "LAURA MILENA"
$asig.Range("H8").Value = "ALEJANDRO"
$asig.Range("H9").Value = "JAIRO ANDRES"
$asig.Range("H10").Value = "BRIAN DAVID"
$asig.Range("H11").Value = "MIGUEL ANGEL"
$asig.Range("H12").Value = "JHONATAN STEVEN"
$asig.Range("H13").Value = "JHOHANNS"

# ------------------------------------------------------------------
# 5. Row 13 closing/footer cell, styled in Courier New
# ------------------------------------------------------------------
$asig.Range("C13").Font.Size = 13
$asig.Range("C13").Font.Color = 0
$asig.Range("C13").Font.Name = "Courier New"
$asig.Rows.Item(13).RowHeight = 18

# ------------------------------------------------------------------
# 6. Stray cell far below, holding a single space character
# ------------------------------------------------------------------
$asig.Range("I22").Value = " "

# ------------------------------------------------------------------
# 7. Column widths / hidden helper columns
# ------------------------------------------------------------------
$asig.Columns.Item("B").ColumnWidth = 14.67
$asig.Columns.Item("C").ColumnWidth = 16.83
$asig.Columns.Item("D").ColumnWidth = 16.83
$asig.Columns.Item("E").ColumnWidth = 10.17
$asig.Columns.Item("G").ColumnWidth = 0
$asig.Columns.Item("G").Hidden = $true
$asig.Columns.Item("H").ColumnWidth = 0
$asig.Columns.Item("H").Hidden = $true

# ------------------------------------------------------------------
# 8. "catalogo" sheet view tweaks - selection moved to I8
# ------------------------------------------------------------------
$catalogo.Range("I8").Select()

# ------------------------------------------------------------------
# 9. Make "Asignacion" the active/selected tab (as in the saved file)
# ------------------------------------------------------------------
$asig.Range("I22").Select()
$asig.Activate()

Write-Output "done"
